$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1728.2307
$ws.Range("I28").Value = 158.5
$ws.Range("K28").Value = 158.5
$ws.Range("M28").Value = 326.5

$ws.Range("H33").Value = 379.9375
$ws.Range("J33").Value = 95
$ws.Range("L33").Value = 95
$ws.Range("N33").Value = -553

$ws.Range("H62").Value = 21816.666
$ws.Range("I62").Value = 7725
$ws.Range("K62").Value = 7725
$ws.Range("M62").Value = -7101

$ws.Range("H65").Value = 21816.666
$ws.Range("I65").Value = 7725
$ws.Range("K65").Value = 38625
$ws.Range("M65").Value = -35505

$ws.Range("H70").Value = 6853.0713
$ws.Range("I70").Value = 1399
$ws.Range("J70").Value = 7272.615
$ws.Range("K70").Value = 4197
$ws.Range("L70").Value = 21817.845
$ws.Range("M70").Value = -3927
$ws.Range("N70").Value = -22357.845

$ws.Range("H73").Value = 6853.0713
$ws.Range("I73").Value = 1399
$ws.Range("J73").Value = 7272.615
$ws.Range("K73").Value = 4197
$ws.Range("L73").Value = 21817.845
$ws.Range("M73").Value = -3261
$ws.Range("N73").Value = -23689.845

$ws.Range("H80").Value = 3987.9546
$ws.Range("I80").Value = 1162.125
$ws.Range("J80").Value = 5602.7144
$ws.Range("K80").Value = 3486.375
$ws.Range("L80").Value = 16808.1432
$ws.Range("M80").Value = -2488.375
$ws.Range("N80").Value = -18804.1432

$ws.Range("H83").Value = 3987.9546
$ws.Range("I83").Value = 1162.125
$ws.Range("J83").Value = 5602.7144
$ws.Range("K83").Value = 10459.125
$ws.Range("L83").Value = 50424.4296
$ws.Range("M83").Value = -5467.125
$ws.Range("N83").Value = -60408.4296

$ws.Range("H96").Value = 514
$ws.Range("I96").Value = 465.4
$ws.Range("K96").Value = 1396.2
$ws.Range("M96").Value = -23.19999999999982

$ws.Range("H106").Value = 166669580
$ws.Range("I106").Value = 166669580
$ws.Range("K106").Value = 166669580
$ws.Range("M106").Value = -166668949

$ws.Range("H112").Value = 5570.654
$ws.Range("J112").Value = 6146
$ws.Range("L112").Value = 18438
$ws.Range("N112").Value = -20654

$ws.Range("H113").Value = 6287.5
$ws.Range("J113").Value = 6602.5
$ws.Range("L113").Value = 6602.5
$ws.Range("N113").Value = -13110.5

$ws.Range("H127").Value = 6246
$ws.Range("I127").Value = 6057.5
$ws.Range("K127").Value = 18172.5
$ws.Range("M127").Value = -13212.5

$ws.Range("H132").Value = 2183.3928
$ws.Range("I132").Value = 2183.3928
$ws.Range("K132").Value = 6550.178400000001
$ws.Range("M132").Value = -4020.178400000001

$ws.Range("H138").Value = 3208.8062
$ws.Range("I138").Value = 1757.3125
$ws.Range("K138").Value = 5271.9375
$ws.Range("M138").Value = -131.9375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2590.61
$ws.Range("I32").Value = 1855.7906
$ws.Range("J32").Value = 7104.5
$ws.Range("K32").Value = 1855.7906
$ws.Range("L32").Value = 7104.5
$ws.Range("M32").Value = -1568.7906
$ws.Range("N32").Value = -7678.5

$ws.Range("H38").Value = 4500
$ws.Range("I38").Value = 4500
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 4500
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -4033
$ws.Range("N38").ClearContents()

$ws.Range("H45").Value = 89069.164
$ws.Range("I45").Value = 168639.67
$ws.Range("J45").Value = 9498.666999999999
$ws.Range("K45").Value = 168639.67
$ws.Range("L45").Value = 9498.666999999999
$ws.Range("M45").Value = -168262.67
$ws.Range("N45").Value = -10252.667

$ws.Range("H61").Value = 1940.5714
$ws.Range("I61").Value = 1940.5714
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1940.5714
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1728.5714
$ws.Range("N61").ClearContents()

$ws.Range("H97").Value = 22651
$ws.Range("I97").Value = 22975.375
$ws.Range("K97").Value = 22975.375
$ws.Range("M97").Value = -22479.375

$ws.Range("H124").Value = 20000
$ws.Range("J124").Value = 20000
$ws.Range("L124").Value = 20000
$ws.Range("N124").Value = -29820

$ws.Range("H132").Value = 3005
$ws.Range("I132").Value = 2273
$ws.Range("K132").Value = 6819
$ws.Range("M132").Value = -4289

$ws.Range("H136").Value = 1940.5714
$ws.Range("I136").Value = 1940.5714
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 5821.7142
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -3271.7142
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 40026.074
$ws.Range("I99").Value = 60936.824
$ws.Range("J99").Value = 4477.8
$ws.Range("K99").Value = 60936.824
$ws.Range("L99").Value = 4477.8
$ws.Range("M99").Value = -59438.824
$ws.Range("N99").Value = -7473.8

$ws.Range("H105").Value = 2222.182
$ws.Range("I105").Value = 1993.7778
$ws.Range("J105").Value = 3250
$ws.Range("K105").Value = 1993.7778
$ws.Range("L105").Value = 3250
$ws.Range("M105").Value = -246.7778000000001
$ws.Range("N105").Value = -6744

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1479.6765
$ws.Range("I58").Value = 1193.8334
$ws.Range("J58").Value = 3623.5
$ws.Range("K58").Value = 1193.8334
$ws.Range("L58").Value = 3623.5
$ws.Range("M58").Value = -990.8334
$ws.Range("N58").Value = -4029.5

$ws.Range("H86").Value = 6703.04
$ws.Range("I86").Value = 6550.857
$ws.Range("J86").Value = 6896.727
$ws.Range("K86").Value = 6550.857
$ws.Range("L86").Value = 6896.727
$ws.Range("M86").Value = -5427.857
$ws.Range("N86").Value = -9142.726999999999

$ws.Range("H89").Value = 6703.04
$ws.Range("I89").Value = 6550.857
$ws.Range("J89").Value = 6896.727
$ws.Range("K89").Value = 32754.285
$ws.Range("L89").Value = 34483.635
$ws.Range("M89").Value = -27138.285
$ws.Range("N89").Value = -45715.635

$ws.Range("H105").Value = 3040.2
$ws.Range("I105").Value = 1876
$ws.Range("K105").Value = 1876
$ws.Range("M105").Value = -129

$ws.Range("H131").Value = 43033.242
$ws.Range("J131").Value = 43033.242
$ws.Range("L131").Value = 43033.242
$ws.Range("N131").Value = -53113.242

$ws.Range("H136").Value = 1479.6765
$ws.Range("I136").Value = 1193.8334
$ws.Range("J136").Value = 3623.5
$ws.Range("K136").Value = 3581.5002
$ws.Range("L136").Value = 10870.5
$ws.Range("M136").Value = -1031.5002
$ws.Range("N136").Value = -15970.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 514998.5
$ws.Range("I62").Value = 514998.5
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 514998.5
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -514312.5
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 514998.5
$ws.Range("I65").Value = 514998.5
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 1544995.5
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -1541563.5
$ws.Range("N65").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1161.75
$ws.Range("I16").Value = 870.2105
$ws.Range("J16").Value = 2269.6
$ws.Range("K16").Value = 870.2105
$ws.Range("L16").Value = 2269.6
$ws.Range("M16").Value = -700.2105
$ws.Range("N16").Value = -2609.6

$ws.Range("H46").Value = 5727.1113
$ws.Range("I46").Value = 4549.1113
$ws.Range("K46").Value = 4549.1113
$ws.Range("M46").Value = -4361.1113

$ws.Range("H100").Value = 53785
$ws.Range("I100").Value = 3971.4285
$ws.Range("K100").Value = 3971.4285
$ws.Range("M100").Value = -3430.4285

$ws.Range("H109").Value = 30964.666
$ws.Range("J109").Value = 30964.666
$ws.Range("L109").Value = 30964.666
$ws.Range("N109").Value = -33738.666

$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws.Range("H122").Value = 5040.8887
$ws.Range("I122").Value = 3760.7058
$ws.Range("K122").Value = 11282.1174
$ws.Range("M122").Value = -8832.117400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H106").Value = 53800
$ws.Range("J106").Value = 53800
$ws.Range("L106").Value = 53800
$ws.Range("N106").Value = -56324

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws.Range("H123").Value = 58714.5
$ws.Range("J123").Value = 58714.5
$ws.Range("L123").Value = 58714.5
$ws.Range("N123").Value = -68514.5

$ws.Range("H132").Value = 275728.53
$ws.Range("I132").Value = 5625.2354
$ws.Range("J132").Value = 3336899.2
$ws.Range("K132").Value = 16875.7062
$ws.Range("L132").Value = 10010697.6
$ws.Range("M132").Value = -14345.7062
$ws.Range("N132").Value = -10015757.6
